# Fix register map: REG_ENCODER1/REG_ENCODER2 are int32 (4 bytes each) but
# were only allocated 2 bytes in the register map, which clobbered
# REG_SPEED1/REG_SPEED2. Re-lay the byte ranges:
#   REG_ENCODER1: 28-29 -> 28-31
#   REG_ENCODER2: 30-31 -> 32-35
#   REG_SPEED1:   32-33 -> 36-37
#   REG_SPEED2:   34-35 -> 38-39

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B19").Value = "28-31"
$ws.Range("B20").Value = "32-35"
$ws.Range("B21").Value = "36-37"
$ws.Range("B22").Value = "38-39"

# Scroll the view so row 13 is at the top, with B23 selected (matches the
# saved view state in the workbook after the edit).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("B23").Select()
